# ClimateScenarios.xml - Correct the paths in scenario #8.
# Scenario.cpp - Correct a logic error ...
# This script reproduces the corresponding spreadsheet edits:
#  - Inserts a new data row (new row 34) with scenario "Baseline 2010 C189"
#    ahead of the existing rows (which all shift down by one).
#  - Appends a new trailing data row (new row 62) with scenario
#    "Demo_Baseline 2010-18 C192".
#  - Updates the frozen-pane / selection view state on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert a new row at position 34. This shifts the current rows
#    34-60 down to 35-61 (Excel automatically keeps their formatting,
#    formulas and relative references correct).
# ---------------------------------------------------------------------
$ws.Rows("34:34").Insert()

# Helper to apply a "highlighted" (yellow fill) numeric style
function Set-NumCell($addr, $value, $fmt, $fill) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $rng.NumberFormat = $fmt
    if ($fill) {
        $rng.Interior.Color = 65535
    } else {
        $rng.Interior.ColorIndex = -4142
    }
}

# ---------------------------------------------------------------------
# New row 34 values
# ---------------------------------------------------------------------
$ws.Range("A34").Value = "CW3M"
$ws.Range("B34").Value = "Baseline 2010 C189"
$ws.Range("C34").Value = 2010

Set-NumCell "D34" 831.51080300000001     "0.00"     $false
Set-NumCell "E34" 1908.5467530000001     "0.00"     $false
Set-NumCell "F34" 1.2276739999999999     "0.00"     $false
Set-NumCell "G34" 302.74935900000003     "0.00"     $true
Set-NumCell "H34" 10.610913999999999     "0.00"     $false
Set-NumCell "I34" 6.4022730000000001     "0.00"     $false
Set-NumCell "J34" 8.8404570000000007     "0.00"     $false
Set-NumCell "K34" 755.04443400000002     "0.00"     $true
Set-NumCell "L34" 59.834083999999997     "0.00"     $false
Set-NumCell "M34" 1296.8793949999999     "0.00"     $true
Set-NumCell "N34" 941.98541299999999     "0.00"     $true
Set-NumCell "O34" 5853.7861329999996     "0"        $true
Set-NumCell "P34" 29450.638672000001     "0"        $false
Set-NumCell "Q34" 1.5360050000000001     "0.00"     $false

$ws.Range("R34").Value = 0.00050199999999999995
$ws.Range("R34").NumberFormat = "General"
$ws.Range("R34").Interior.ColorIndex = -4142

$ws.Range("S34").Value = 2010

# ---------------------------------------------------------------------
# 2. Append a brand-new row 62 (after the shifted data, which now ends
#    at row 61) with scenario "Demo_Baseline 2010-18 C192".
# ---------------------------------------------------------------------
$ws.Range("A62").Value = "CW3M"
$ws.Range("B62").Value = "Demo_Baseline 2010-18 C192"
$ws.Range("C62").Value = "2010-18"

Set-NumCell "D62" 1000.3124864444443     "0.00"     $true
Set-NumCell "E62" 1763.5263265555557     "0.00"     $false
Set-NumCell "F62" 0.999942                "0.00"     $false
Set-NumCell "G62" 305.6782124444444      "0.00"     $true
Set-NumCell "H62" 9.775355222222224      "0.00"     $false
Set-NumCell "I62" 6.8948233333333331     "0.00"     $false
Set-NumCell "J62" 8.145128999999999      "0.00"     $false
Set-NumCell "K62" 673.17452677777771     "0.00"     $true
Set-NumCell "L62" 60.018756111111117     "0.00"     $false
Set-NumCell "M62" 1321.9402533333332     "0.00"     $true
Set-NumCell "N62" 1024.1975572222223     "0.00"     $true
Set-NumCell "O62" 4583.9874403333333     "0"        $true
Set-NumCell "P62" 27227.338324888889     "0"        $false
Set-NumCell "Q62" 0.28907633333333327    "0.00"     $false
Set-NumCell "R62" 0.00006222222222222222 "0.000000" $false

$ws.Range("S62").Value = "2010-18"

# ---------------------------------------------------------------------
# 3. Update the sheet view: frozen pane top-left cell and the active
#    selection, to match the new layout/size of the data.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("A50").Select()
$win.Panes.Item($win.Panes.Count).ScrollRow = 50
$ws.Range("S63").Select()

$wb.Save()
